$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.040.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.414.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("E11").Value = "  -3.18%  "
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.06%  "
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.854.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.977.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.410.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "321.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("E25").Value = "  -5.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "568.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0932"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("E34").Value = "  -2.38%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.380"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "148.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0531"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.594"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0920"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("E51").Value = "  +0.69%  "
